# Factor the "group" classifier column out of the formula column, on the
# first two sheets: insert two new columns (C and D) that both hold the
# literal group name, shift the old "=F/G" lookup formula two columns to
# the right, and fill in the previously-missing row-10 counters.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Group classification per data row (rows 2..15), same on both sheets.
$groups = @{
    2  = "group3"
    3  = "group1"
    4  = "group3"
    5  = "group3"
    6  = "group1"
    7  = "group1"
    8  = "group2"
    9  = "group2"
    10 = "group2"
    11 = "group2"
    12 = "group2"
    13 = "group2"
    14 = "group2"
    15 = "group1"
}

# --- Sheet 1 ("sheet1") ----------------------------------------------------
# Insert two blank columns before column C; this pushes the old "=F{r}"
# formula (and its value column) two columns to the right, becoming
# "=H{r}" with the values now living in column H.
$ws1.Range("C:D").Insert() | Out-Null

foreach ($r in 2..15) {
    $g = $groups[$r]
    $ws1.Cells.Item($r, 3).Value = $g   # C{r}
    $ws1.Cells.Item($r, 4).Value = $g   # D{r}
}

# The inserted columns pick up the (highlighted) row format on the
# duplicate-name rows 14/15; re-apply the plain formula-column formatting
# (column E, which is what column C used to look like) so C/D don't stay
# highlighted red like A/B on those rows.
$ws1.Range("E14:E15").Copy() | Out-Null
$ws1.Range("C14:D15").PasteSpecial(-4122) | Out-Null
$ws1.Application.CutCopyMode = $false

# Row 10 was missing its counter value.
$ws1.Range("B10").Value = 9

# --- Sheet 2 ("Feuille2") ---------------------------------------------------
# Row 10 was missing its counter value and its group classification.
$ws2.Range("B10").Value = 9
$ws2.Range("C10").Value = "group2"

# Move the selection on sheet 2 to C11, then restore sheet 1 as the active
# (selected) sheet/tab.
$ws2.Range("C11").Select() | Out-Null
$ws1.Activate() | Out-Null
